# Update the 3-1_DataTable worksheet with the latest APS (Annual Population
# Survey) data: "Latest period" moves from Jan-Dec 2024 to Apr 2024-Mar 2025,
# and "Next period" moves from Apr 2024-Mar 2025 to Jul 2024-Jun 2025.
# This affects the three rows that use the APS data source (Employment
# volumes, Employment by occupation, Employment by industry - rows 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$latestPeriod = "Apr 2024 - Mar 2025 (17/07/25)"
$nextPeriod = "Jul 2024 - Jun 2025 (14/10/25)"

$ws.Range("C2").Value = $latestPeriod
$ws.Range("D2").Value = $nextPeriod

$ws.Range("C3").Value = $latestPeriod
$ws.Range("D3").Value = $nextPeriod

$ws.Range("C4").Value = $latestPeriod
$ws.Range("D4").Value = $nextPeriod

# Update the saved selection/view to cell C2 (also clears the scrolled
# topLeftCell that pointed at row 4).
$ws.Range("C2").Select()
